# ---------------------------------------------------------------------------
# C5-PowerPoint.pptx edit:
#   1. The table on slide 6 switches from the custom table style
#      {86F2E2C2-2B39-4DAA-8BE7-44ADCEC6F6C6} (defined in ppt/tableStyles.xml)
#      to the built-in gallery style {CE37C3BB-766A-4281-88D5-652027EBA842}.
#   2. The deck's theme colour scheme changes from the "Integral" palette to
#      the default "Office" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-style the table on slide 6 (it is the only table in the deck).
# ---------------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{CE37C3BB-766A-4281-88D5-652027EBA842}", $true)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the theme palette back to the standard Office colours.
#    ColorScheme.Colors(1..12) map 1:1 onto dk1, lt1, dk2, lt2, accent1-6,
#    hlink, folHlink in the theme part backing the slide master.
#    RGB values are packed 0x00BBGGRR (standard VBA/COM RGB()).
# ---------------------------------------------------------------------------
$officeThemeColors = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

$master = $p.SlideMaster
$colorScheme = $master.ColorScheme
for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
